$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.019.57'
$ws.Range("E2").Value = '  -0.63%  '
$ws.Range("D3").Value = '1.828.59'
$ws.Range("E3").Value = '  -0.45%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.99'
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6227'
$ws.Range("E6").Value = '  -5.78%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.80'
$ws.Range("E8").Value = '  +7.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07538'
$ws.Range("E9").Value = '  +1.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.2904'
$ws.Range("E10").Value = '  -0.91%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.70'
$ws.Range("E11").Value = '  -1.38%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07641'
$ws.Range("E12").Value = '  -1.47%  '
$ws.Range("D13").Value = '1.836.49'
$ws.Range("E13").Value = '  -13.69%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.950'
$ws.Range("E14").Value = '  -0.79%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6636'
$ws.Range("E15").Value = '  -0.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.20'
$ws.Range("E16").Value = '  -0.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000009089'
$ws.Range("E17").Value = '  +5.93%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.966'
$ws.Range("E18").Value = '  -2.69%  '
$ws.Range("D19").Value = '28.930.90'
$ws.Range("E19").Value = '  -0.83%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '224.30'
$ws.Range("E20").Value = '  -1.40%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.31'
$ws.Range("E21").Value = '  -1.47%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.180'
$ws.Range("E23").Value = '  +0.69%  '
$ws.Range("E24").Value = '  +0.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '159.76'
$ws.Range("E25").Value = '  +0.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.386'
$ws.Range("E26").Value = '  -2.66%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1355'
$ws.Range("E27").Value = '  -3.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.79'
$ws.Range("E28").Value = '  -1.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.493'
$ws.Range("E29").Value = '  -1.57%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.020'
$ws.Range("E30").Value = '  -0.81%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.037'
$ws.Range("E31").Value = '  -2.05%  '
$ws.Range("E32").Value = '  +0.52%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05189'
$ws.Range("E33").Value = '  -1.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.839'
$ws.Range("E34").Value = '  -1.74%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.149'
$ws.Range("E35").Value = '  +0.20%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7285'
$ws.Range("E36").Value = '  -0.96%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.614'
$ws.Range("E37").Value = '  -1.45%  '
$ws.Range("D38").Value = '1.276.53'
$ws.Range("E38").Value = '  -2.36%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.759'
$ws.Range("E39").Value = '  +0.79%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01785'
$ws.Range("E40").Value = '  -0.74%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.390'
$ws.Range("E41").Value = '  +5.62%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8897'
$ws.Range("E42").Value = '  -3.76%  '
$ws.Range("E43").Value = '  +0.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.33'
$ws.Range("E44").Value = '  -1.17%  '
$ws.Range("D45").Value = '1.979.27'
$ws.Range("E45").Value = '  +1.72%  '
$ws.Range("E46").Value = '  -0.74%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '63.36'
$ws.Range("E47").Value = '  -0.80%  '
$ws.Range("E48").Value = '  -3.67%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07431'
$ws.Range("E49").Value = '  -16.23%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3971'
$ws.Range("E50").Value = '  -1.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.844'
$ws.Range("E51").Value = '  +1.19%  '
